{"js": "// Office.js (Word JavaScript API) script\n// Applies the LOM3241.docx edit described by the diff:\n//  1. Heading3 title: \"Chemistry of Materials\" -> \"Materials chemistry\"\n//  2. \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\"\n//  3. Insert an italic English translation paragraph after the \"Objetivos\" paragraph\n//  4. Replace the \"Programa resumido\" paragraph text and insert an italic English translation after it\n//  5. Replace the \"Programa\" paragraph text (collapsing the <w:br/> separated runs into one run)\n//     and insert an italic English translation after it\n//  6. Update the \"Crit\u00e9rio\" evaluation text\n//  7. Update the \"Norma de recupera\u00e7\u00e3o\" text\n\nconst T = JSON.parse(`{\"HEADING3_OLD\": \"Chemistry of Materials\", \"HEADING3_NEW\": \"Materials chemistry\", \"ATIVACAO_OLD\": \"Ativa\u00e7\u00e3o: 01/01/2012\", \"ATIVACAO_NEW\": \"Ativa\u00e7\u00e3o: 01/01/2023\", \"OBJETIVOS_PT\": \"Fornecer ao estudante os principais tipos de s\u00edntese org\u00e2nica e inorg\u00e2nica de materiais bem como apresentar as principais t\u00e9cnicas anal\u00edticas para caracteriza\u00e7\u00e3o de materiais.\", \"OBJETIVOS_EN\": \"Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization.\", \"RESUMIDO_OLD\": \"Qu\u00edmica e materiais. Liga\u00e7\u00f5es qu\u00edmicas. S\u00edntese de materiais e transforma\u00e7\u00f5es qu\u00edmicas. Processos. T\u00e9cnicas de caracteriza\u00e7\u00e3o de materiais. Tipos de materiais. Considera\u00e7\u00f5es econ\u00f4micas e ambientais.\", \"RESUMIDO_NEW_PT\": \"Introdu\u00e7\u00e3o \u00e0 qu\u00edmica e sua associa\u00e7\u00e3o com s\u00edntese de novos materiais. A vis\u00e3o moderna do \u00e1tomo  e Liga\u00e7\u00f5es qu\u00edmicas. Estrutura cristalina e t\u00e9cnicas de caracteriza\u00e7\u00e3o cristalogr\u00e1fica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, s\u00edntese e aplica\u00e7\u00f5es. Processos e T\u00e9cnicas de crescimento de cristais de um modo geral. Pol\u00edmeros condutores e suas aplica\u00e7\u00f5es em tecnologica moderna.\", \"RESUMIDO_NEW_EN\": \"Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology.\", \"PROGRAMA_NEW_PT\": \"Qu\u00edmica de materiais: defini\u00e7\u00e3o; papel da qu\u00edmica na ci\u00eancia de materiais; fundamentos.Atom\u00edstica e a vis\u00e3o moderna do \u00e1tomo com fundamentos qu\u00e2nticos.Tipos de liga\u00e7\u00f5es qu\u00edmicas: for\u00e7as de van der Waals, potencial de Lennard-Jones, liga\u00e7\u00e3o covalente, liga\u00e7\u00f5es por coordena\u00e7\u00e3o, liga\u00e7\u00f5es i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.Materiais policristalinos e monocristalinos. A ordem cristalogr\u00e1fica e t\u00e9cnicas de caracteriza\u00e7\u00e3o cristalogr\u00e1fica e microsc\u00f3pica. A import\u00e2ncia de monocristais em aplica\u00e7\u00f5es eletr\u00f4nicas. T\u00e9cnicas de crescimento de cristais de alta qualidade tais como: m\u00e9todo do fluxo, m\u00e9todo Czochralski, m\u00e9todo Brigdmann, m\u00e9todo do transporte de vapor e m\u00e9todo de crescimento de transporte de vapor modificado e isot\u00e9rmico. Materiais amorfos e sua import\u00e2ncia para a tecnologica moderna. Conceitos e t\u00e9cnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, t\u00e9cnicas de crescimento tais como: vapor qu\u00edmico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletr\u00f3lise para revestimento protetivo, conceitos e aplica\u00e7\u00f5es. S\u00edntese de pol\u00edmeros condutores, conceitos e aplica\u00e7\u00f5es como dispositivos eletr\u00f4nicos.\", \"PROGRAMA_NEW_EN\": \"Materials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices.\", \"CRITERIO_OLD\": \"M\u00e9dia ponderada de duas provas escritas, trabalhos e relat\u00f3rios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4\", \"CRITERIO_NEW\": \"M\u00e9dia simples de duas provas escritas,  Conceito Final = (P1 + P2)/2\", \"RECUPERACAO_OLD\": \"Aplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\", \"RECUPERACAO_NEW\": \"Aplica\u00e7\u00e3o de duas provas escritas dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo.\"}`);\n\n// ---- 1. Heading3 title ------------------------------------------------\nconst headingResults = context.document.body.search(T.HEADING3_OLD, { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\nif (headingResults.items.length === 0) {\n  throw new Error(\"Could not find heading text to replace: \" + T.HEADING3_OLD);\n}\nheadingResults.items[0].insertText(T.HEADING3_NEW, Word.InsertLocation.replace);\n\n// ---- 2. Ativa\u00e7\u00e3o date ---------------------------------------------------\nconst ativacaoResults = context.document.body.search(T.ATIVACAO_OLD, { matchCase: true });\nativacaoResults.load(\"text\");\nawait context.sync();\nif (ativacaoResults.items.length === 0) {\n  throw new Error(\"Could not find Ativa\u00e7\u00e3o text to replace: \" + T.ATIVACAO_OLD);\n}\nativacaoResults.items[0].insertText(T.ATIVACAO_NEW, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---- 3. Objetivos: add italic English translation paragraph ------------\nconst objetivosResults = context.document.body.search(T.OBJETIVOS_PT, { matchCase: true });\nobjetivosResults.load(\"text\");\nawait context.sync();\nif (objetivosResults.items.length === 0) {\n  throw new Error(\"Could not find Objetivos paragraph text: \" + T.OBJETIVOS_PT);\n}\nconst objetivosPara = objetivosResults.items[0].paragraphs.getFirst();\nconst objetivosEnPara = objetivosPara.insertParagraph(T.OBJETIVOS_EN, Word.InsertLocation.after);\nobjetivosEnPara.font.set({ italic: true });\nawait context.sync();\n\n// ---- 4. Programa resumido: replace text + add italic translation -------\nconst resumidoResults = context.document.body.search(T.RESUMIDO_OLD, { matchCase: true });\nresumidoResults.load(\"text\");\nawait context.sync();\nif (resumidoResults.items.length === 0) {\n  throw new Error(\"Could not find Programa resumido paragraph text: \" + T.RESUMIDO_OLD);\n}\nconst resumidoRange = resumidoResults.items[0];\nconst resumidoPara = resumidoRange.paragraphs.getFirst();\nresumidoRange.insertText(T.RESUMIDO_NEW_PT, Word.InsertLocation.replace);\nawait context.sync();\nconst resumidoEnPara = resumidoPara.insertParagraph(T.RESUMIDO_NEW_EN, Word.InsertLocation.after);\nresumidoEnPara.font.set({ italic: true });\nawait context.sync();\n\n// ---- 5. Programa: replace text (merge into a single run) + translation -\n// Locate the \"Programa\" heading paragraph, then its following content paragraph.\nconst headings = context.document.body.paragraphs;\nheadings.load(\"text,style\");\nawait context.sync();\nlet programaBodyPara = null;\nfor (let i = 0; i < headings.items.length; i++) {\n  if (headings.items[i].style === \"Heading 2\" && headings.items[i].text === \"Programa\") {\n    programaBodyPara = headings.items[i + 1];\n    break;\n  }\n}\nif (!programaBodyPara) {\n  throw new Error(\"Could not find the 'Programa' section body paragraph\");\n}\nprogramaBodyPara.insertText(T.PROGRAMA_NEW_PT, Word.InsertLocation.replace);\nawait context.sync();\nconst programaEnPara = programaBodyPara.insertParagraph(T.PROGRAMA_NEW_EN, Word.InsertLocation.after);\nprogramaEnPara.font.set({ italic: true });\nawait context.sync();\n\n// ---- 6. Crit\u00e9rio: update evaluation text --------------------------------\nconst criterioResults = context.document.body.search(T.CRITERIO_OLD, { matchCase: true });\ncriterioResults.load(\"text\");\nawait context.sync();\nif (criterioResults.items.length === 0) {\n  throw new Error(\"Could not find Crit\u00e9rio text to replace: \" + T.CRITERIO_OLD);\n}\ncriterioResults.items[0].insertText(T.CRITERIO_NEW, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---- 7. Norma de recupera\u00e7\u00e3o: update text -------------------------------\nconst recuperacaoResults = context.document.body.search(T.RECUPERACAO_OLD, { matchCase: true });\nrecuperacaoResults.load(\"text\");\nawait context.sync();\nif (recuperacaoResults.items.length === 0) {\n  throw new Error(\"Could not find Norma de recupera\u00e7\u00e3o text to replace: \" + T.RECUPERACAO_OLD);\n}\nrecuperacaoResults.items[0].insertText(T.RECUPERACAO_NEW, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# PowerShell (Word COM interop) script\n# Applies the LOM3241.docx edit described by the diff:\n#  1. Heading3 title: \"Chemistry of Materials\" -> \"Materials chemistry\"\n#  2. \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\"\n#  3. Insert an italic English translation paragraph after the \"Objetivos\" paragraph\n#  4. Replace the \"Programa resumido\" paragraph text and insert an italic English translation after it\n#  5. Replace the \"Programa\" paragraph text (collapsing the <w:br/> separated runs into one run)\n#     and insert an italic English translation after it\n#  6. Update the \"Crit\u00e9rio\" evaluation text\n#  7. Update the \"Norma de recupera\u00e7\u00e3o\" text\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText {\n    param(\n        $Doc,\n        [string]$OldText,\n        [string]$NewText\n    )\n    $find = $Doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $OldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $NewText\n    $ok = $find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)\n    if (-not $ok) {\n        throw \"Replace-ExactText: could not find text: $OldText\"\n    }\n}\n\nfunction Add-ItalicTranslationAfter {\n    # Finds the paragraph whose trimmed text equals $MatchText, replaces its text with\n    # $NewText (if not $null), then inserts a new paragraph right after it containing\n    # $TranslationText formatted in italics (run-level only, no paragraph-mark formatting).\n    param(\n        $Doc,\n        [string]$MatchText,\n        [string]$NewText,\n        [string]$TranslationText\n    )\n    $count = $Doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $Doc.Paragraphs($i)\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $MatchText) {\n            if ($NewText -ne $null) {\n                $p.Range.Text = $NewText\n            }\n            $p.Range.InsertParagraphAfter() | Out-Null\n            $newPara = $Doc.Paragraphs($i + 1)\n            $newPara.Range.Text = $TranslationText\n            $italicRange = $Doc.Range($newPara.Range.Start, $newPara.Range.End - 1)\n            $italicRange.Font.Italic = 1\n            return\n        }\n    }\n    throw \"Add-ItalicTranslationAfter: could not find paragraph matching: $MatchText\"\n}\n\n# ---- 1. Heading3 title ---------------------------------------------------\nReplace-ExactText $d @'\nChemistry of Materials\n'@ @'\nMaterials chemistry\n'@\n\n# ---- 2. Ativa\u00e7\u00e3o date -----------------------------------------------------\nReplace-ExactText $d @'\nAtiva\u00e7\u00e3o: 01/01/2012\n'@ @'\nAtiva\u00e7\u00e3o: 01/01/2023\n'@\n\n# ---- 3. Objetivos: add italic English translation paragraph ---------------\nAdd-ItalicTranslationAfter $d @'\nFornecer ao estudante os principais tipos de s\u00edntese org\u00e2nica e inorg\u00e2nica de materiais bem como apresentar as principais t\u00e9cnicas anal\u00edticas para caracteriza\u00e7\u00e3o de materiais.\n'@ $null @'\nProvide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization.\n'@\n\n# ---- 4. Programa resumido: replace text + add italic translation ----------\nAdd-ItalicTranslationAfter $d @'\nQu\u00edmica e materiais. Liga\u00e7\u00f5es qu\u00edmicas. S\u00edntese de materiais e transforma\u00e7\u00f5es qu\u00edmicas. Processos. T\u00e9cnicas de caracteriza\u00e7\u00e3o de materiais. Tipos de materiais. Considera\u00e7\u00f5es econ\u00f4micas e ambientais.\n'@ @'\nIntrodu\u00e7\u00e3o \u00e0 qu\u00edmica e sua associa\u00e7\u00e3o com s\u00edntese de novos materiais. A vis\u00e3o moderna do \u00e1tomo  e Liga\u00e7\u00f5es qu\u00edmicas. Estrutura cristalina e t\u00e9cnicas de caracteriza\u00e7\u00e3o cristalogr\u00e1fica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, s\u00edntese e aplica\u00e7\u00f5es. Processos e T\u00e9cnicas de crescimento de cristais de um modo geral. Pol\u00edmeros condutores e suas aplica\u00e7\u00f5es em tecnologica moderna.\n'@ @'\nIntroduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology.\n'@\n\n# ---- 5. Programa: replace text (merge into a single run) + translation ----\n$programaOld = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Programa\") {\n        $programaOld = $d.Paragraphs($i + 1).Range.Text.TrimEnd([char]13, [char]7)\n        break\n    }\n}\nif ($programaOld -eq $null) {\n    throw \"Could not find the 'Programa' section body paragraph\"\n}\nAdd-ItalicTranslationAfter $d $programaOld @'\nQu\u00edmica de materiais: defini\u00e7\u00e3o; papel da qu\u00edmica na ci\u00eancia de materiais; fundamentos.Atom\u00edstica e a vis\u00e3o moderna do \u00e1tomo com fundamentos qu\u00e2nticos.Tipos de liga\u00e7\u00f5es qu\u00edmicas: for\u00e7as de van der Waals, potencial de Lennard-Jones, liga\u00e7\u00e3o covalente, liga\u00e7\u00f5es por coordena\u00e7\u00e3o, liga\u00e7\u00f5es i\u00f4nicas e liga\u00e7\u00f5es met\u00e1licas.Materiais policristalinos e monocristalinos. A ordem cristalogr\u00e1fica e t\u00e9cnicas de caracteriza\u00e7\u00e3o cristalogr\u00e1fica e microsc\u00f3pica. A import\u00e2ncia de monocristais em aplica\u00e7\u00f5es eletr\u00f4nicas. T\u00e9cnicas de crescimento de cristais de alta qualidade tais como: m\u00e9todo do fluxo, m\u00e9todo Czochralski, m\u00e9todo Brigdmann, m\u00e9todo do transporte de vapor e m\u00e9todo de crescimento de transporte de vapor modificado e isot\u00e9rmico. Materiais amorfos e sua import\u00e2ncia para a tecnologica moderna. Conceitos e t\u00e9cnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, t\u00e9cnicas de crescimento tais como: vapor qu\u00edmico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletr\u00f3lise para revestimento protetivo, conceitos e aplica\u00e7\u00f5es. S\u00edntese de pol\u00edmeros condutores, conceitos e aplica\u00e7\u00f5es como dispositivos eletr\u00f4nicos.\n'@ @'\nMaterials chemistry: definition; role of chemistry in materials science; fundamentals.Atomistics and the modern view of the atom with quantum foundations.Types of chemical bonds: van der Waals forces, Lennard-Jones potential, covalent bonding, coordination bonds, ionic bonds and metallic bonds.Polycrystalline and monocrystalline materials. The crystallographic order and crystallographic and microscopic characterization techniques. The importance of single crystals in electronic applications. High quality crystal growth techniques such as: flow method, Czochralski method, Brigdmann method, vapor transport method and modified isothermal vapor transport growth method. Amorphous materials and their importance for modern technology. Concepts and techniques for growing amorphous materials. Epitaxial thin films, growth techniques such as: chemical vapor, sputtering, laser ablation and MBE. Thin films grown by electrolysis for protective coating, concepts and applications. Synthesis of conductive polymers, concepts and applications as electronic devices.\n'@\n\n# ---- 6. Crit\u00e9rio: update evaluation text -----------------------------------\nReplace-ExactText $d @'\nM\u00e9dia ponderada de duas provas escritas, trabalhos e relat\u00f3rios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4\n'@ @'\nM\u00e9dia simples de duas provas escritas,  Conceito Final = (P1 + P2)/2\n'@\n\n# ---- 7. Norma de recupera\u00e7\u00e3o: update text ----------------------------------\nReplace-ExactText $d @'\nAplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\n'@ @'\nAplica\u00e7\u00e3o de duas provas escritas dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo.\n'@\n"}
